# Insert a new row at position 156; this shifts existing rows 156:254
# down to 157:255 (Excel also extends the used range / dimension
# automatically, and copies formatting from the row being pushed down).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("156:156").Insert()

# Populate the newly inserted row 156 with the new record.
$ws.Cells.Item(156, 1).Value  = 7
$ws.Cells.Item(156, 2).Value  = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(156, 3).Value  = "Ñuble"
$ws.Cells.Item(156, 4).Value  = 44606
$ws.Cells.Item(156, 5).Value  = 16
$ws.Cells.Item(156, 6).Value  = 100114013
$ws.Cells.Item(156, 7).Value  = "Zanahoria"
$ws.Cells.Item(156, 8).Value  = "Sin especificar"
$ws.Cells.Item(156, 9).Value  = "Primera"
$ws.Cells.Item(156, 10).Value = 100
$ws.Cells.Item(156, 11).Value = 6000
$ws.Cells.Item(156, 12).Value = 6500
$ws.Cells.Item(156, 13).Value = 6250
$ws.Cells.Item(156, 14).Value = "$/saco 20 kilos"
$ws.Cells.Item(156, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(156, 16).Value = 312
$ws.Cells.Item(156, 17).Value = 20
$ws.Cells.Item(156, 18).Value = "Hortaliza"

# Make sure the date cell keeps the date number format used by the
# rest of column D.
$ws.Cells.Item(156, 4).NumberFormat = $ws.Cells.Item(157, 4).NumberFormat
